$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 55 (item id 5517)
$ws.Range("H55").Value = 330.7143
$ws.Range("I55").Value = 260.08334
$ws.Range("J55").Value = 383.6875
$ws.Range("K55").Value = 260.08334
$ws.Range("L55").Value = 383.6875
$ws.Range("M55").Value = -46.08334000000002
$ws.Range("N55").Value = -811.6875
# row 112 (item id 27960)
$ws.Range("H112").Value = 1371.6552
$ws.Range("J112").Value = 1402.7858
$ws.Range("L112").Value = 4208.357400000001
$ws.Range("N112").Value = -6424.357400000001
# row 116 (item id 27778)
$ws.Range("H116").Value = 593620.5
$ws.Range("I116").Value = 911049.8
$ws.Range("K116").Value = 911049.8
$ws.Range("M116").Value = -907607.8
# row 141 (item id 44161)
$ws.Range("H141").Value = 52409.05
$ws.Range("I141").Value = 78837
$ws.Range("K141").Value = 236511
$ws.Range("M141").Value = -231331

$ws = $wb.Worksheets.Item("ARM")
# row 61 (item id 43999)
$ws.Range("H61").Value = 1146.6471
$ws.Range("I61").Value = 1093.3125
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1093.3125
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -881.3125
$ws.Range("N61").Value = -2424
# row 136 (item id 43999)
$ws.Range("H136").Value = 1146.6471
$ws.Range("I136").Value = 1093.3125
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3279.9375
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -729.9375
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
# row 53 (item id 27158)
$ws.Range("H53").Value = 42666.668
$ws.Range("J53").Value = 42666.668
$ws.Range("L53").Value = 42666.668
$ws.Range("N53").Value = -43814.668
# row 95 (item id 18194)
$ws.Range("H95").Value = 32500
$ws.Range("J95").Value = 32500
$ws.Range("L95").Value = 32500
$ws.Range("N95").Value = -37992
# row 134 (item id 43998)
$ws.Range("H134").Value = 2618.5527
$ws.Range("I134").Value = 1315.9166
$ws.Range("J134").Value = 4851.643
$ws.Range("K134").Value = 3947.7498
$ws.Range("L134").Value = 14554.929
$ws.Range("M134").Value = -1412.7498
$ws.Range("N134").Value = -19624.929

$ws = $wb.Worksheets.Item("CRP")
# row 31 (item id 44023)
$ws.Range("H31").Value = 153488.34
$ws.Range("I31").Value = 356274.66
$ws.Range("J31").Value = 2392.6667
$ws.Range("K31").Value = 356274.66
$ws.Range("L31").Value = 2392.6667
$ws.Range("M31").Value = -355979.66
$ws.Range("N31").Value = -2982.6667
# row 34 (item id 44023)
$ws.Range("H34").Value = 153488.34
$ws.Range("I34").Value = 356274.66
$ws.Range("J34").Value = 2392.6667
$ws.Range("K34").Value = 356274.66
$ws.Range("L34").Value = 2392.6667
$ws.Range("M34").Value = -356072.66
$ws.Range("N34").Value = -2796.6667
# row 134 (item id 44020)
$ws.Range("H134").Value = 1686.4231
$ws.Range("I134").Value = 1029.2354
$ws.Range("J134").Value = 2927.7778
$ws.Range("K134").Value = 3087.7062
$ws.Range("L134").Value = 8783.3334
$ws.Range("M134").Value = -552.7062000000001
$ws.Range("N134").Value = -13853.3334

$ws = $wb.Worksheets.Item("CUL")
# row 68 (item id 12895)
$ws.Range("H68").Value = 1050.5807
$ws.Range("I68").Value = 768.5599999999999
$ws.Range("K68").Value = 2305.68
$ws.Range("M68").Value = -1494.68
# row 71 (item id 12895)
$ws.Range("H71").Value = 1050.5807
$ws.Range("I71").Value = 768.5599999999999
$ws.Range("K71").Value = 6917.039999999999
$ws.Range("M71").Value = -2861.039999999999
# row 113 (item id 27843)
$ws.Range("H113").Value = 448.73077
$ws.Range("I113").Value = 469.4
$ws.Range("J113").Value = 429.5926
$ws.Range("K113").Value = 1408.2
$ws.Range("L113").Value = 1288.7778
$ws.Range("M113").Value = 761.8000000000002
$ws.Range("N113").Value = -5628.7778
# row 119 (item id 27873)
$ws.Range("H119").Value = 1250
$ws.Range("I119").Value = 1250
$ws.Range("K119").Value = 3750
$ws.Range("M119").Value = 1088
# row 121 (item id 27878)
$ws.Range("H121").Value = 1680.8379
$ws.Range("J121").Value = 1907.1774
$ws.Range("L121").Value = 5721.5322
$ws.Range("N121").Value = -8341.5322

$ws = $wb.Worksheets.Item("GSM")
# row 18 (item id 4309)
$ws.Range("H18").Value = 22300
$ws.Range("J18").Value = 22300
$ws.Range("L18").Value = 22300
$ws.Range("N18").Value = -22886
# row 126 (item id 36184)
$ws.Range("H126").Value = 4094.28
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 6039.6665
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 18118.9995
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -23058.9995
# row 140 (item id 42458)
$ws.Range("H140").Value = 39302.855
$ws.Range("J140").Value = 39302.855
$ws.Range("L140").Value = 39302.855
$ws.Range("N140").Value = -49662.855

$ws = $wb.Worksheets.Item("LTW")
# row 29 (item id 3576)
$ws.Range("H29").Value = 24999.5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 24999.5
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 24999.5
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -25589.5
# row 33 (item id 4106)
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
# row 61 (item id 27740)
$ws.Range("H61").Value = 2222.5
$ws.Range("I61").Value = 2222.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2222.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2020.5
$ws.Range("N61").ClearContents()
# row 64 (item id 10810)
$ws.Range("H64").Value = 32300
$ws.Range("J64").Value = 32300
$ws.Range("L64").Value = 32300
$ws.Range("N64").Value = -32750
# row 67 (item id 10810)
$ws.Range("H67").Value = 32300
$ws.Range("J67").Value = 32300
$ws.Range("L67").Value = 32300
$ws.Range("N67").Value = -33860
# row 94 (item id 18067)
$ws.Range("H94").Value = 32483.334
$ws.Range("J94").Value = 32483.334
$ws.Range("L94").Value = 32483.334
$ws.Range("N94").Value = -33835.334
# row 113 (item id 27740)
$ws.Range("H113").Value = 2222.5
$ws.Range("I113").Value = 2222.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2222.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -52.5
$ws.Range("N113").ClearContents()
# row 132 (item id 44058)
$ws.Range("H132").Value = 3149.7188
$ws.Range("I132").Value = 2494.1
$ws.Range("J132").Value = 4242.4165
$ws.Range("K132").Value = 7482.299999999999
$ws.Range("L132").Value = 12727.2495
$ws.Range("M132").Value = -4952.299999999999
$ws.Range("N132").Value = -17787.2495
# row 136 (item id 44060)
$ws.Range("H136").Value = 3876.9644
$ws.Range("I136").Value = 1276.1538
$ws.Range("K136").Value = 3828.4614
$ws.Range("M136").Value = -1278.4614

$ws = $wb.Worksheets.Item("WVR")
# row 126 (item id 36210)
$ws.Range("H126").Value = 711659.4399999999
$ws.Range("I126").Value = 1674.5834
$ws.Range("J126").Value = 3551599
$ws.Range("K126").Value = 5023.7502
$ws.Range("L126").Value = 10654797
$ws.Range("M126").Value = -2553.7502
$ws.Range("N126").Value = -10659737
# row 132 (item id 44029)
$ws.Range("H132").Value = 7411975
$ws.Range("I132").Value = 6306.25
$ws.Range("J132").Value = 13336510
$ws.Range("K132").Value = 18918.75
$ws.Range("L132").Value = 40009530
$ws.Range("M132").Value = -16388.75
$ws.Range("N132").Value = -40014590
# row 136 (item id 44031)
$ws.Range("H136").Value = 9145.096
$ws.Range("I136").Value = 9675.333000000001
$ws.Range("J136").Value = 8438.111000000001
$ws.Range("K136").Value = 29025.999
$ws.Range("L136").Value = 25314.333
$ws.Range("M136").Value = -26475.999
$ws.Range("N136").Value = -30414.333
